# Applies the "updated cryptos list" refresh described by the commit.
# For each changed row we update Coin (B) / Link (C) when the rank swapped,
# and always refresh Price (D) and Volume(1h) (E). Price values that look like
# plain numbers are forced back to literal text (NumberFormat "@") so Excel
# does not silently coerce them into numeric cells, then the style is reset
# to "Normal" so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText($row, $text) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$rows = @(
    @{ Row = 2; D = '64.049.48'; E = '  +1.42%  ' },
    @{ Row = 3; D = '3.320.44'; E = '  +5.97%  ' },
    @{ Row = 4; E = '  -0.04%  ' },
    @{ Row = 5; D = '599.65'; E = '  +1.01%  ' },
    @{ Row = 6; D = '143.80'; E = '  +5.26%  ' },
    @{ Row = 7; E = '  -0.04%  ' },
    @{ Row = 8; D = '3.321.85'; E = '  +6.25%  ' },
    @{ Row = 9; D = '0.525'; E = '  +1.39%  ' },
    @{ Row = 10; E = '  +3.03%  ' },
    @{ Row = 11; D = '5.53'; E = '  +5.85%  ' },
    @{ Row = 12; D = '0.477'; E = '  +4.22%  ' },
    @{ Row = 13; D = '0.0000251'; E = '  +1.43%  ' },
    @{ Row = 14; D = '34.91'; E = '  +1.82%  ' },
    @{ Row = 15; D = '3.864.97'; E = '  +6.28%  ' },
    @{ Row = 16; E = '  +1.15%  ' },
    @{ Row = 17; D = '3.316.88'; E = '  +6.24%  ' },
    @{ Row = 18; D = '64.085.65'; E = '  +1.56%  ' },
    @{ Row = 19; D = '6.93'; E = '  +3.30%  ' },
    @{ Row = 20; D = '484.49'; E = '  +2.08%  ' },
    @{ Row = 21; D = '14.32'; E = '  +0.29%  ' },
    @{ Row = 22; E = '  +5.93%  ' },
    @{ Row = 23; D = '8.04'; E = '  +3.86%  ' },
    @{ Row = 24; D = '13.66'; E = '  +4.65%  ' },
    @{ Row = 25; D = '84.86'; E = '  -2.38%  ' },
    @{ Row = 27; D = '2.79'; E = '  +2.35%  ' },
    @{ Row = 28; E = '  +1.44%  ' },
    @{ Row = 29; B = 'RenderToken'; C = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D = '8.27'; E = '  +4.08%  ' },
    @{ Row = 30; B = 'FirstDigitalUSD'; C = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'; D = '1.00'; E = '  -0.04%  ' },
    @{ Row = 31; D = '29.59'; E = '  +10.54%  ' },
    @{ Row = 32; E = '  +5.44%  ' },
    @{ Row = 33; D = '0.107'; E = '  -0.59%  ' },
    @{ Row = 34; D = '2.57'; E = '  +1.55%  ' },
    @{ Row = 35; E = '  +2.95%  ' },
    @{ Row = 36; D = '6.03'; E = '  +3.33%  ' },
    @{ Row = 37; D = '53.39'; E = '  +2.55%  ' },
    @{ Row = 38; D = '0.0₃0761'; E = '  +7.82%  ' },
    @{ Row = 39; E = '  +3.64%  ' },
    @{ Row = 40; D = '433.96'; E = '  +2.86%  ' },
    @{ Row = 41; D = '3.045.79'; E = '  +5.45%  ' },
    @{ Row = 42; B = 'dogwifhat'; C = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'; D = '2.80'; E = '  +3.73%  ' },
    @{ Row = 43; B = 'Cosmos'; C = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D = '8.48'; E = '  +2.61%  ' },
    @{ Row = 44; E = '  -2.36%  ' },
    @{ Row = 45; D = '0.269'; E = '  +2.12%  ' },
    @{ Row = 46; D = '2.23'; E = '  +4.50%  ' },
    @{ Row = 47; D = '26.69'; E = '  +3.31%  ' },
    @{ Row = 48; B = 'Arweave'; C = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'; D = '36.01'; E = '  +12.52%  ' },
    @{ Row = 49; B = 'USDe'; C = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'; D = '0.999'; E = '  +0.03%  ' },
    @{ Row = 50; B = 'Stellar'; C = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D = '0.116'; E = '  +2.47%  ' },
    @{ Row = 51; D = '2.33'; E = '  +1.82%  ' }
)

foreach ($r in $rows) {
    if ($r.ContainsKey("B")) { $ws.Cells.Item($r.Row, 2).Value = $r.B }
    if ($r.ContainsKey("C")) { $ws.Cells.Item($r.Row, 3).Value = $r.C }
    if ($r.ContainsKey("D")) { Set-PriceText $r.Row $r.D }
    if ($r.ContainsKey("E")) { $ws.Cells.Item($r.Row, 5).Value = $r.E }
}
